$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SignIn")

# Update the hyperlink target for the "Url" cell (A2) to the new local test URL,
# keeping it as the same hyperlink object (just a new Address), then update the
# visible cell text to match.
foreach ($h in $ws.Hyperlinks) {
    if ($h.Range.Address() -eq "`$A`$2") {
        $h.Address = "http://localhost:5000/"
    }
}

$ws.Range("A2").Value = "http://localhost:5000/"
$ws.Range("B2").Value = "aswini.sanal@outlook.com"
$ws.Range("C2").Value = "Testing0123*"

# SignIn becomes the active/selected sheet.
$ws.Activate()
